$wb = $excel.ActiveWorkbook

# --- Transactions sheet: set its selection before leaving it, so the
#     selection persisted in the file matches the target (D5), and it will
#     no longer be the active/tabSelected sheet once we move on.
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Activate() | Out-Null
$wsTrans.Range("D5").Select() | Out-Null

# --- Repayment schedule sheet: insert a new blank column before the old
#     "Late" column (old N), shifting Late/heading(4)/Outstanding from
#     N:P to O:Q.
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate() | Out-Null
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 9.9

# Final selection on this sheet becomes the active one (activeTab moves to
# this sheet, tabSelected flips here too).
$ws.Range("S9").Select() | Out-Null
